$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D8').Value = '12 ماهه منتهی به 1397/12'
$ws.Range('E8').Value = '12 ماهه منتهی به 1398/12'
$ws.Range('F8').Value = '12 ماهه منتهی به 1399/12'
$ws.Range('G8').Value = '12 ماهه منتهی به 1400/12'
$ws.Range('H8').Value = '12 ماهه منتهی به 1401/12'
$ws.Range('D9').Value = '1399-04-16 (8)'
$ws.Range('E9').Value = '1400-04-20 (8)'
$ws.Range('F9').Value = '1401-04-11 (8)'
$ws.Range('G9').Value = '1401-10-28 (6)'
$ws.Range('H9').NumberFormat = "@"
$ws.Range('H9').Value = '1402-02-28'
$ws.Range('D12').Value = 192976
$ws.Range('E12').Value = 823607
$ws.Range('F12').Value = 585849
$ws.Range('G12').Value = 737383
$ws.Range('H12').Value = 2195096
$ws.Range('D13').Value = 0
$ws.Range('E13').Value = 0
$ws.Range('F13').Value = 0
$ws.Range('G13').Value = 0
$ws.Range('H13').Value = 0
$ws.Range('D14').Value = 226351
$ws.Range('E14').Value = 151741
$ws.Range('F14').Value = 200851
$ws.Range('G14').Value = 406773
$ws.Range('H14').Value = 876880
$ws.Range('D15').Value = 890953
$ws.Range('E15').Value = 1683851
$ws.Range('F15').Value = 3820584
$ws.Range('G15').Value = 4296489
$ws.Range('H15').Value = 5009042
$ws.Range('D16').Value = 448500
$ws.Range('E16').Value = 178750
$ws.Range('F16').Value = 990998
$ws.Range('G16').Value = 608423
$ws.Range('H16').Value = 865856
$ws.Range('D17').Value = 0
$ws.Range('E17').Value = 0
$ws.Range('F17').Value = 0
$ws.Range('G17').Value = 0
$ws.Range('H17').Value = 0
$ws.Range('D18').Value = 1758780
$ws.Range('E18').Value = 2837949
$ws.Range('F18').Value = 5598282
$ws.Range('G18').Value = 6049068
$ws.Range('H18').Value = 8946874
$ws.Range('D19').Value = 0
$ws.Range('E19').Value = 0
$ws.Range('F19').Value = 0
$ws.Range('G19').Value = 0
$ws.Range('H19').Value = 0
$ws.Range('D20').Value = 13014
$ws.Range('E20').Value = 16823
$ws.Range('F20').Value = 139567
$ws.Range('G20').Value = 47999
$ws.Range('H20').Value = 46219
$ws.Range('D21').Value = 0
$ws.Range('E21').Value = 0
$ws.Range('F21').Value = 0
$ws.Range('G21').Value = 0
$ws.Range('H21').Value = 0
$ws.Range('D22').Value = 336019
$ws.Range('E22').Value = 629991
$ws.Range('F22').Value = 12019491
$ws.Range('G22').Value = 13267752
$ws.Range('H22').Value = 15069287
$ws.Range('D23').Value = 11431
$ws.Range('E23').Value = 18758
$ws.Range('F23').Value = 16586
$ws.Range('G23').Value = 58552
$ws.Range('H23').Value = 58470
$ws.Range('D24').Value = '-'
$ws.Range('E24').Value = '-'
$ws.Range('F24').Value = '-'
$ws.Range('G24').Value = '-'
$ws.Range('H24').Value = '-'
$ws.Range('D25').Value = 0
$ws.Range('E25').Value = 0
$ws.Range('F25').Value = 0
$ws.Range('G25').Value = 0
$ws.Range('H25').Value = 0
$ws.Range('D26').Value = 360464
$ws.Range('E26').Value = 665572
$ws.Range('F26').Value = 12175644
$ws.Range('G26').Value = 13374303
$ws.Range('H26').Value = 15173976
$ws.Range('D27').Value = 2119244
$ws.Range('E27').Value = 3503521
$ws.Range('F27').Value = 17773926
$ws.Range('G27').Value = 19423371
$ws.Range('H27').Value = 24120850
$ws.Range('D29').Value = 742227
$ws.Range('E29').Value = 1003058
$ws.Range('F29').Value = 1222000
$ws.Range('G29').Value = 2052023
$ws.Range('H29').Value = 3294892
$ws.Range('D30').Value = '-'
$ws.Range('E30').Value = '-'
$ws.Range('F30').Value = '-'
$ws.Range('G30').Value = '-'
$ws.Range('H30').Value = '-'
$ws.Range('D31').Value = 411370
$ws.Range('E31').Value = 886059
$ws.Range('F31').Value = 924852
$ws.Range('G31').Value = 642513
$ws.Range('H31').Value = 1338241
$ws.Range('D32').Value = 106880
$ws.Range('E32').Value = 276982
$ws.Range('F32').Value = 475401
$ws.Range('G32').Value = 347584
$ws.Range('H32').Value = 387307
$ws.Range('D33').Value = 16026
$ws.Range('E33').Value = 25755
$ws.Range('F33').Value = 77059
$ws.Range('G33').Value = 152813
$ws.Range('H33').Value = 172782
$ws.Range('D34').Value = 242007
$ws.Range('E34').Value = 57243
$ws.Range('F34').Value = 273853
$ws.Range('G34').Value = 2113531
$ws.Range('H34').Value = 2075847
$ws.Range('D35').Value = 0
$ws.Range('E35').Value = 0
$ws.Range('F35').Value = 0
$ws.Range('G35').Value = 0
$ws.Range('H35').Value = 0
$ws.Range('D36').Value = 0
$ws.Range('E36').Value = 0
$ws.Range('F36').Value = 0
$ws.Range('G36').Value = 0
$ws.Range('H36').Value = 0
$ws.Range('D37').Value = 1518510
$ws.Range('E37').Value = 2249097
$ws.Range('F37').Value = 2973165
$ws.Range('G37').Value = 5308464
$ws.Range('H37').Value = 7269069
$ws.Range('D38').Value = 0
$ws.Range('E38').Value = 0
$ws.Range('F38').Value = 1219595
$ws.Range('G38').Value = 0
$ws.Range('H38').Value = 0
$ws.Range('D39').Value = '-'
$ws.Range('E39').Value = '-'
$ws.Range('F39').Value = '-'
$ws.Range('G39').Value = '-'
$ws.Range('H39').Value = '-'
$ws.Range('D40').Value = 0
$ws.Range('E40').Value = 0
$ws.Range('F40').Value = 0
$ws.Range('G40').Value = 0
$ws.Range('H40').Value = 2000000
$ws.Range('D41').Value = 65956
$ws.Range('E41').Value = 79895
$ws.Range('F41').Value = 110295
$ws.Range('G41').Value = 143358
$ws.Range('H41').Value = 203052
$ws.Range('D42').Value = 65956
$ws.Range('E42').Value = 79895
$ws.Range('F42').Value = 1329890
$ws.Range('G42').Value = 143358
$ws.Range('H42').Value = 2203052
$ws.Range('D43').Value = 1584466
$ws.Range('E43').Value = 2328992
$ws.Range('F43').Value = 4303055
$ws.Range('G43').Value = 5451822
$ws.Range('H43').Value = 9472121
$ws.Range('D45').Value = 286469
$ws.Range('E45').Value = 286469
$ws.Range('F45').Value = 9882730
$ws.Range('G45').Value = 9882730
$ws.Range('H45').Value = 13000000
$ws.Range('D46').Value = 0
$ws.Range('E46').Value = 0
$ws.Range('F46').Value = 0
$ws.Range('G46').Value = 0
$ws.Range('H46').Value = 0
$ws.Range('D47').Value = 0
$ws.Range('E47').Value = 0
$ws.Range('F47').Value = 0
$ws.Range('G47').Value = 0
$ws.Range('H47').Value = 0
$ws.Range('D48').Value = 0
$ws.Range('E48').Value = 0
$ws.Range('F48').Value = -349921
$ws.Range('G48').Value = -595305
$ws.Range('H48').Value = -595305
$ws.Range('D49').Value = 0
$ws.Range('E49').Value = 0
$ws.Range('F49').Value = 111056
$ws.Range('G49').Value = 25796
$ws.Range('H49').Value = 25796
$ws.Range('D50').Value = 22576
$ws.Range('E50').Value = 28647
$ws.Range('F50').Value = 205530
$ws.Range('G50').Value = 266862
$ws.Range('H50').Value = 339721
$ws.Range('D51').Value = 0
$ws.Range('E51').Value = 0
$ws.Range('F51').Value = 0
$ws.Range('G51').Value = 0
$ws.Range('H51').Value = 0
$ws.Range('D52').Value = '-'
$ws.Range('E52').Value = '-'
$ws.Range('F52').Value = '-'
$ws.Range('G52').Value = '-'
$ws.Range('H52').Value = '-'
$ws.Range('D53').Value = 0
$ws.Range('E53').Value = 0
$ws.Range('F53').Value = 0
$ws.Range('G53').Value = 0
$ws.Range('H53').Value = 0
$ws.Range('D54').Value = '-'
$ws.Range('E54').Value = '-'
$ws.Range('F54').Value = '-'
$ws.Range('G54').Value = '-'
$ws.Range('H54').Value = '-'
$ws.Range('D55').Value = 0
$ws.Range('E55').Value = 0
$ws.Range('F55').Value = 0
$ws.Range('G55').Value = 0
$ws.Range('H55').Value = 0
$ws.Range('D56').Value = 225733
$ws.Range('E56').Value = 859413
$ws.Range('F56').Value = 3621476
$ws.Range('G56').Value = 4391466
$ws.Range('H56').Value = 1878517
$ws.Range('D57').Value = 534778
$ws.Range('E57').Value = 1174529
$ws.Range('F57').Value = 13470871
$ws.Range('G57').Value = 13971549
$ws.Range('H57').Value = 14648729
$ws.Range('D58').Value = 2119244
$ws.Range('E58').Value = 3503521
$ws.Range('F58').Value = 17773926
$ws.Range('G58').Value = 19423371
$ws.Range('H58').Value = 24120850
